$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.929.93'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '2.208.92'
$ws.Range('E3').Value = '  -1.62%  '
$ws.Range('D5').Value = "'230.20"
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').Value = "'0.616"
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  -1.93%  '
$ws.Range('D7').Value = "'60.53"
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  -1.14%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').Value = "'0.401"
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  -1.21%  '
$ws.Range('D10').Value = "'0.0897"
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +1.82%  '
$ws.Range('E11').Value = '  -0.66%  '
$ws.Range('D12').Value = '2.543.51'
$ws.Range('E12').Value = '  -1.54%  '
$ws.Range('D13').Value = "'15.39"
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  -2.92%  '
$ws.Range('D14').Value = "'22.00"
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('D15').Value = "'0.796"
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').Value = "'5.56"
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  -0.85%  '
$ws.Range('D17').Value = '2.213.72'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').Value = '41.865.69'
$ws.Range('E18').Value = '  +0.19%  '
$ws.Range('D19').Value = '0.0₃0938'
$ws.Range('E19').Value = '  +4.38%  '
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('D21').Value = "'6.05"
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  -1.75%  '
$ws.Range('D22').Value = "'242.12"
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  -2.98%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  -0.15%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = "'2.38"
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +0.10%  '
$ws.Range('B25').Value = 'PancakeSwap'
$ws.Range('C25').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D25').Value = "'2.35"
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  -1.75%  '
$ws.Range('D26').Value = "'9.61"
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  -0.39%  '
$ws.Range('D27').Value = "'168.66"
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('E28').Value = '  -1.38%  '
$ws.Range('D29').Value = "'20.40"
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +1.06%  '
$ws.Range('D30').Value = "'1.42"
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  -1.24%  '
$ws.Range('E31').Value = '  -5.31%  '
$ws.Range('E32').Value = '  -2.42%  '
$ws.Range('D33').Value = "'4.94"
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  -4.13%  '
$ws.Range('D34').Value = "'4.59"
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('D35').Value = "'0.0645"
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +2.74%  '
$ws.Range('D36').Value = "'6.26"
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  -7.01%  '
$ws.Range('D37').Value = "'3.52"
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  -7.45%  '
$ws.Range('D38').Value = "'2.32"
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  -3.26%  '
$ws.Range('B39').Value = 'BinanceUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D39').Value = "'1.00"
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = "'0.0243"
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +2.11%  '
$ws.Range('D41').Value = "'0.000228"
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  -9.37%  '
$ws.Range('D42').Value = "'8.52"
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  -2.75%  '
$ws.Range('D43').Value = "'0.0950"
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  -4.17%  '
$ws.Range('E44').Value = '  -0.47%  '
$ws.Range('B45').Value = 'FTXToken'
$ws.Range('C45').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D45').Value = "'4.38"
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  -12.40%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').Value = "'96.27"
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -4.41%  '
$ws.Range('D47').Value = '1.455.03'
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('D48').Value = "'2.74"
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  -1.52%  '
$ws.Range('D49').Value = "'15.98"
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  -3.66%  '
$ws.Range('E50').Value = '  -3.48%  '
$ws.Range('D51').Value = "'2.20"
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +1.83%  '
